$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three changed values (trim trailing digit/precision) while
# leaving the rest of the "C Lag"/"A Lag" table untouched.
$ws.Range("B3").Value = "-2.82***"
$ws.Range("C2").Value = "-0.01*"
$ws.Range("C3").Value = "-0.47***"
